$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.256.11'
$ws.Range("E2").Value = '  -1.41%  '
$ws.Range("D3").Value = '3.620.32'
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.46'
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.74'
$ws.Range("E6").Value = '  +2.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.613'
$ws.Range("E7").Value = '  -2.41%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.673'
$ws.Range("E9").Value = '  -5.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.145'
$ws.Range("E10").Value = '  -9.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.78'
$ws.Range("E11").Value = '  -3.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000255'
$ws.Range("E12").Value = '  -11.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.95'
$ws.Range("E13").Value = '  -5.81%  '
$ws.Range("D14").Value = '4.177.71'
$ws.Range("E14").Value = '  -1.62%  '
$ws.Range("D15").Value = '3.600.14'
$ws.Range("E15").Value = '  -1.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.126'
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.41'
$ws.Range("E17").Value = '  -4.23%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '67.000.15'
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.25'
$ws.Range("E19").Value = '  -3.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.06'
$ws.Range("E20").Value = '  -4.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '392.72'
$ws.Range("E21").Value = '  -3.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.31'
$ws.Range("E22").Value = '  -4.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '84.89'
$ws.Range("E23").Value = '  -3.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.87'
$ws.Range("E24").Value = '  -4.54%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.24'
$ws.Range("E25").Value = '  -3.26%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.05'
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.34'
$ws.Range("E27").Value = '  -3.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.62'
$ws.Range("E28").Value = '  -11.53%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.96'
$ws.Range("E29").Value = '  -4.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.18'
$ws.Range("E30").Value = '  -4.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.75'
$ws.Range("E31").Value = '  -5.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '65.38'
$ws.Range("E32").Value = '  +1.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.89'
$ws.Range("E33").Value = '  -3.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '593.23'
$ws.Range("E34").Value = '  +0.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.113'
$ws.Range("E35").Value = '  -3.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '41.47'
$ws.Range("E36").Value = '  -3.43%  '
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.375'
$ws.Range("E39").Value = '  -5.40%  '
$ws.Range("D40").Value = '0.0₃0746'
$ws.Range("E40").Value = '  -15.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.132'
$ws.Range("E41").Value = '  -2.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.78'
$ws.Range("E42").Value = '  -7.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0411'
$ws.Range("E43").Value = '  -5.31%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.41'
$ws.Range("E44").Value = '  -10.61%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.685.15'
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.130'
$ws.Range("E46").Value = '  -3.04%  '
$ws.Range("E47").Value = '  -3.92%  '
$ws.Range("E48").Value = '  -6.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '136.06'
$ws.Range("E49").Value = '  -3.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.24'
$ws.Range("E50").Value = '  -8.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.58'
$ws.Range("E51").Value = '  -5.98%  '
